# Apply Configurable zero_before_threshold parameter updates to the
# Step3_DataPts_* sheets (0.5 / 0.7 / 0.8 / 0.9).
#
# Only columns C (First_Noticeable_Increase_Index),
# E (First_Noticeable_Increase_Cumulative_Value) and
# G (Pulse_Width) change, for rows 3-6 on each sheet.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> row -> column letter -> new value
$updates = @{
    "Step3_DataPts_0.5" = @{
        3 = @{ C = 47; E = 0.0008117423756870443; G = 48 }
        4 = @{ C = 52; E = 0.05869535002919633;   G = 46 }
        5 = @{ C = 48; E = 0.007701964088017953;  G = 49 }
        6 = @{ C = 48; E = 0.007985624599653402;  G = 49 }
    }
    "Step3_DataPts_0.7" = @{
        3 = @{ C = 47; E = 0.0008117423756870443; G = 61 }
        4 = @{ C = 52; E = 0.05869535002919633;   G = 61 }
        5 = @{ C = 48; E = 0.007701964088017953;  G = 61 }
        6 = @{ C = 48; E = 0.007985624599653402;  G = 61 }
    }
    "Step3_DataPts_0.8" = @{
        3 = @{ C = 47; E = 0.0008117423756870443; G = 76 }
        4 = @{ C = 52; E = 0.05869535002919633;   G = 73 }
        5 = @{ C = 48; E = 0.007701964088017953;  G = 76 }
        6 = @{ C = 48; E = 0.007985624599653402;  G = 77 }
    }
    "Step3_DataPts_0.9" = @{
        3 = @{ C = 47; E = 0.0008117423756870443; G = 112 }
        4 = @{ C = 52; E = 0.05869535002919633;   G = 109 }
        5 = @{ C = 48; E = 0.007701964088017953;  G = 114 }
        6 = @{ C = 48; E = 0.007985624599653402;  G = 115 }
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $cols = $rows[$row]
        foreach ($col in $cols.Keys) {
            $ws.Range("$col$row").Value = $cols[$col]
        }
    }
}
